$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old sample data rows (2 and 3) - the sheet becomes header-only
# ---------------------------------------------------------------------------
$ws.Rows("2:3").Delete()

# ---------------------------------------------------------------------------
# 2. Rewrite the header row (row 1) with the new set of questions
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "ФИО"
$ws.Range("D1").Value = "Телефон"
$ws.Range("E1").Value = "Гос. Знак"
$ws.Range("F1").Value = "Грузоподьемность"
$ws.Range("G1").Value = "Измерения"
$ws.Range("H1").Value = "Кузов"
$ws.Range("I1").Value = "Город"
$ws.Range("J1").Value = "Дистанция"
$ws.Range("K1").Value = "ЮР. Статус"
$ws.Range("L1").Value = "Владение "
$ws.Range("M1").Value = "Тип загрузки"
$ws.Range("N1").Value = "Грузы"

# ---------------------------------------------------------------------------
# 3. New header cells (J1:N1) get the same bold/centered look as the rest of
#    row 1 (copy format from A1), but with a left+right only thin border
#    instead of the full box used by A1:I1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("J1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($col in @("J", "K", "L", "M", "N")) {
    $cell = $ws.Range($col + "1")
    $cell.Borders.Item(8).LineStyle = -4142
    $cell.Borders.Item(9).LineStyle = -4142
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------------
# 4. Column widths for the new/changed columns
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 14.166666666666666
$ws.Columns("D").ColumnWidth = 13.053385416666666
$ws.Columns("E").ColumnWidth = 17.276041666666668
$ws.Columns("F").ColumnWidth = 19.166666666666668
$ws.Columns("G").ColumnWidth = 14.276041666666666
$ws.Columns("I").ColumnWidth = 10.276041666666666
$ws.Columns("J").ColumnWidth = 13.944010416666666
$ws.Columns("K").ColumnWidth = 12.608072916666666
$ws.Columns("L").ColumnWidth = 15.721354166666666
$ws.Columns("M").ColumnWidth = 11.385416666666666

# ---------------------------------------------------------------------------
# 5. Selection / active cell, matching the updated worksheet view
# ---------------------------------------------------------------------------
$ws.Range("A2:L2").Select()

Write-Output "done"
